{"js": "// Update the date line and the 25 multiplication problems (three-digit x\n// one-digit) to the new values from the commit.  Every \"old\" string below\n// is unique within the document, and none of the \"new\" strings collide\n// with any \"old\" string, so a straightforward literal search & replace of\n// each pair is safe and order-independent.\nconst replacements = [\n  [\"2025-03-22 Saturday\", \"2025-03-23 Sunday\"],\n  [\"624\u00d79=\", \"466\u00d78=\"],\n  [\"523\u00d77=\", \"609\u00d73=\"],\n  [\"550\u00d77=\", \"411\u00d75=\"],\n  [\"466\u00d77=\", \"550\u00d78=\"],\n  [\"452\u00d72=\", \"959\u00d79=\"],\n  [\"238\u00d72=\", \"609\u00d76=\"],\n  [\"732\u00d78=\", \"945\u00d73=\"],\n  [\"486\u00d75=\", \"907\u00d74=\"],\n  [\"260\u00d75=\", \"143\u00d78=\"],\n  [\"278\u00d79=\", \"838\u00d73=\"],\n  [\"744\u00d74=\", \"563\u00d75=\"],\n  [\"822\u00d76=\", \"195\u00d76=\"],\n  [\"273\u00d77=\", \"943\u00d77=\"],\n  [\"539\u00d78=\", \"383\u00d73=\"],\n  [\"796\u00d78=\", \"524\u00d75=\"],\n  [\"843\u00d72=\", \"558\u00d77=\"],\n  [\"357\u00d75=\", \"396\u00d74=\"],\n  [\"852\u00d76=\", \"396\u00d75=\"],\n  [\"464\u00d75=\", \"273\u00d78=\"],\n  [\"727\u00d74=\", \"660\u00d72=\"],\n  [\"698\u00d76=\", \"241\u00d75=\"],\n  [\"526\u00d77=\", \"424\u00d72=\"],\n  [\"725\u00d75=\", \"136\u00d76=\"],\n  [\"993\u00d76=\", \"970\u00d76=\"],\n  [\"748\u00d72=\", \"483\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems (three-digit x\n# one-digit) to the new values from the commit. Every \"old\" string is\n# unique within the document, and none of the \"new\" strings collide with\n# any \"old\" string, so Find/Replace-All of each literal pair, run in any\n# order, reproduces the target content exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"2025-03-22 Saturday\"; New = \"2025-03-23 Sunday\" },\n  @{ Old = \"624\u00d79=\"; New = \"466\u00d78=\" },\n  @{ Old = \"523\u00d77=\"; New = \"609\u00d73=\" },\n  @{ Old = \"550\u00d77=\"; New = \"411\u00d75=\" },\n  @{ Old = \"466\u00d77=\"; New = \"550\u00d78=\" },\n  @{ Old = \"452\u00d72=\"; New = \"959\u00d79=\" },\n  @{ Old = \"238\u00d72=\"; New = \"609\u00d76=\" },\n  @{ Old = \"732\u00d78=\"; New = \"945\u00d73=\" },\n  @{ Old = \"486\u00d75=\"; New = \"907\u00d74=\" },\n  @{ Old = \"260\u00d75=\"; New = \"143\u00d78=\" },\n  @{ Old = \"278\u00d79=\"; New = \"838\u00d73=\" },\n  @{ Old = \"744\u00d74=\"; New = \"563\u00d75=\" },\n  @{ Old = \"822\u00d76=\"; New = \"195\u00d76=\" },\n  @{ Old = \"273\u00d77=\"; New = \"943\u00d77=\" },\n  @{ Old = \"539\u00d78=\"; New = \"383\u00d73=\" },\n  @{ Old = \"796\u00d78=\"; New = \"524\u00d75=\" },\n  @{ Old = \"843\u00d72=\"; New = \"558\u00d77=\" },\n  @{ Old = \"357\u00d75=\"; New = \"396\u00d74=\" },\n  @{ Old = \"852\u00d76=\"; New = \"396\u00d75=\" },\n  @{ Old = \"464\u00d75=\"; New = \"273\u00d78=\" },\n  @{ Old = \"727\u00d74=\"; New = \"660\u00d72=\" },\n  @{ Old = \"698\u00d76=\"; New = \"241\u00d75=\" },\n  @{ Old = \"526\u00d77=\"; New = \"424\u00d72=\" },\n  @{ Old = \"725\u00d75=\"; New = \"136\u00d76=\" },\n  @{ Old = \"993\u00d76=\"; New = \"970\u00d76=\" },\n  @{ Old = \"748\u00d72=\"; New = \"483\u00d73=\" }\n)\n\nforeach ($r in $replacements) {\n  $find = $d.Content.Find\n  $find.Text = $r.Old\n  $find.Replacement.Text = $r.New\n  $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
